$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added as the 6th data point for this
# subset. It belongs right after the current row 81 (date 44174) and
# before the old row 82 (date 44596), so insert a fresh row at 82 and
# push the rest (old rows 82-86) down to 83-87.
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new record.
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value = 44769
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = 100112031
$ws.Range("G82").Value = "Poroto verde"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 30
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = 30000
$ws.Range("N82").Value = "`$/malla 25 kilos"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 1200
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
